$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column O, shifting the Zone $/TB table (O:S) to (P:T)
$ws.Columns("O:O").Insert()

# Set header and data for the new "Measured [MB/s]" column (N)
$ws.Range("N4").Value = "Measured [MB/s]"
$ws.Range("N5").Value = $ws.Range("C5").Value
$ws.Range("N6").Value = $ws.Range("C6").Value
$ws.Range("N7").Value = $ws.Range("C7").Value
$ws.Range("N8").Value = $ws.Range("C8").Value
$ws.Range("N9").Value = $ws.Range("C9").Value
$ws.Range("N10").Value = $ws.Range("C10").Value
$ws.Range("N11").Value = $ws.Range("C11").Value
$ws.Range("N12").Value = $ws.Range("C12").Value
$ws.Range("N13").Value = $ws.Range("C13").Value
$ws.Range("N14").Value = $ws.Range("C14").Value

$ws.Columns("N").ColumnWidth = 18.28515625
$ws.Columns("T").ColumnWidth = 17
